# Investment_Cost_overview.xlsx update
# - Extend Table1 with a new "Notes" column (A1:H24 -> A1:I24)
# - Fill in investment-cost data for 8 new technologies (rows 14-21)
# - Add Notes text for a few rows, and blank (but styled) Notes cells for the rest
# - Adjust column I width / wrap formatting
# - Update the view selection on sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. Grow the table one column to the right and name the new column "Notes"
# ---------------------------------------------------------------------------
$tbl.Resize($ws.Range("A1:I24"))
$ws.Range("I1").Value = "Notes"

# ---------------------------------------------------------------------------
# 2. New technology rows (14-21)
# ---------------------------------------------------------------------------

# Row 14 - Fischer_Tropsch_Unit
$ws.Range("A14").Value = "Fischer_Tropsch_Unit"
$ws.Range("B14").Value = 2100000
$ws.Range("C14").Value = 1850000
$ws.Range("D14").Value = 1600000
$ws.Range("E14").Value = 1100000
$ws.Range("F14").Style = "Normal"
$ws.Range("F14").Value = 900000
$ws.Range("F14").VerticalAlignment = -4160
$ws.Range("F14").WrapText = $true
$ws.Range("G14").Value = "25Y"
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "No value given for 2025 (linearly approximated)"

# Row 15 - RWGS_Unit
$ws.Range("A15").Value = "RWGS_Unit"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "25Y"
$ws.Range("H15").Value = 0

# Row 16 - Distillation_tower_FT
$ws.Range("A16").Value = "Distillation_tower_FT"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = "25Y"
$ws.Range("H16").Value = 0

# Row 17 - ASU
$ws.Range("A17").Value = "ASU"
$ws.Range("B17").Value = 413513.5
$ws.Range("C17").Value = 377027
$ws.Range("D17").Value = 340540.5
$ws.Range("E17").Value = 267567.59999999998
$ws.Range("F17").Value = 218567.6
$ws.Range("G17").Value = "30Y"
$ws.Range("H17").Value = 0.03
$ws.Range("I17").Value = "Inv_costs calculated using percentage from catalogue (see data_needed)"

# Row 18 - Haber_Bosch_Reactor
$ws.Range("A18").Value = "Haber_Bosch_Reactor"
$ws.Range("B18").Value = 1700000
$ws.Range("C18").Value = 1550000
$ws.Range("D18").Value = 1400000
$ws.Range("E18").Value = 1100000
$ws.Range("F18").Value = 900000
$ws.Range("G18").Value = "30Y"
$ws.Range("H18").Value = 0.03
$ws.Range("I18").Value = "No value given for 2025 (linearly approximated)"

# Row 19 - Methane_Plant
$ws.Range("A19").Value = "Methane_Plant"
$ws.Range("B19").Value = 9600000
$ws.Range("C19").Value = 2100000
$ws.Range("D19").Value = 1500000
$ws.Range("E19").Value = 1200000
$ws.Range("F19").Value = 1200000
$ws.Range("G19").Value = "25Y"
$ws.Range("H19").Value = 0.03
$ws.Range("I19").Value = "Percentage not specifically given"

# Row 20 - Digester
$ws.Range("A20").Value = "Digester"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = "25Y"
$ws.Range("H20").Value = 0.03

# Row 21 - CO2_Remover
$ws.Range("A21").Value = "CO2_Remover"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = "25Y"
$ws.Range("H21").Value = 0.03

# ---------------------------------------------------------------------------
# 3. Format the (mostly empty) Notes column I2:I24 like the rest of the table
#    (same "#,##0.00" number style used throughout columns B:F)
# ---------------------------------------------------------------------------
$ws.Range("I2:I24").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------------
# 4. Column width for the Notes column
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 44

# ---------------------------------------------------------------------------
# 5. Sheet view: scroll so column D is left-most and select J16
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("J16").Select()
